$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.175534725189209
$ws.Range("B1").Value = 4.552795886993408
$ws.Range("C1").Value = 3.511415243148804
$ws.Range("D1").Value = 3.249753475189209
$ws.Range("E1").Value = 2.544822931289673
